# Apply the changes described by the diff:
#  - Insert 3 new "NOT MAPPED" lookup rows into "SAM CPHHolding" (DISEASE_TYPE, INTERVAL, INTERVAL_UNIT_OF_TIME)
#    before the existing row 30 (HoldingStartDate), shifting subsequent rows down.
#  - Update various sheet view (scroll/selection) states across sheets.
#  - Update the workbook window size/position.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1. SAM CPHHolding - insert 3 new rows with "NOT MAPPED" formatting (style
#    copied from the existing template row 23 which has the same pattern).
# ---------------------------------------------------------------------------
$wsHolding = $wb.Worksheets.Item("SAM CPHHolding")

$wsHolding.Rows("30:32").Insert()

# Copy formatting (style 11 for A/C columns, style 12 for E/F columns) from
# the existing "NOT MAPPED" template row (row 23) onto the 3 new rows.
$wsHolding.Range("A23:F23").Copy()
$wsHolding.Range("A30:F30").PasteSpecial(-4122)
$wsHolding.Range("A23:F23").Copy()
$wsHolding.Range("A31:F31").PasteSpecial(-4122)
$wsHolding.Range("A23:F23").Copy()
$wsHolding.Range("A32:F32").PasteSpecial(-4122)

# These template rows only use columns A, C, E and F - clear the spurious
# B/D cells that Insert/PasteSpecial created.
$wsHolding.Range("B30").Clear()
$wsHolding.Range("D30").Clear()
$wsHolding.Range("B31").Clear()
$wsHolding.Range("D31").Clear()
$wsHolding.Range("B32").Clear()
$wsHolding.Range("D32").Clear()

# Set the new cell values.
$wsHolding.Range("A30").Value = "DISEASE_TYPE"
$wsHolding.Range("C30").Value = "NOT MAPPED"

$wsHolding.Range("A31").Value = "INTERVAL"
$wsHolding.Range("C31").Value = "NOT MAPPED"

$wsHolding.Range("A32").Value = "INTERVAL_UNIT_OF_TIME"
$wsHolding.Range("C32").Value = "NOT MAPPED"

# Update this sheet's view (scroll position / active selection).
$wsHolding.Application.ActiveWindow.TabSelected = 1
$wsHolding.Range("A16").Select()
$wsHolding.Application.ActiveWindow.ScrollRow = 16
$wsHolding.Range("A33").Select()

# ---------------------------------------------------------------------------
# 2. SAM Herd - selection moves from C15 to C16.
# ---------------------------------------------------------------------------
$wsHerd = $wb.Worksheets.Item("SAM Herd")
$wsHerd.Range("C16").Select()

# ---------------------------------------------------------------------------
# 3. SAM CPHHolder - scroll position changes (topLeftCell becomes A13).
# ---------------------------------------------------------------------------
$wsCPHHolder = $wb.Worksheets.Item("SAM CPHHolder")
$wsCPHHolder.Range("C23").Select()
$wsCPHHolder.Application.ActiveWindow.ScrollRow = 13

# ---------------------------------------------------------------------------
# 4. SAM Party - scroll resets to top, selection moves from C10 to B27.
# ---------------------------------------------------------------------------
$wsParty = $wb.Worksheets.Item("SAM Party")
$wsParty.Application.ActiveWindow.ScrollRow = 1
$wsParty.Range("B27").Select()

# ---------------------------------------------------------------------------
# 5. Party Roles - no longer the tab-selected sheet.
# ---------------------------------------------------------------------------
$wsPartyRoles = $wb.Worksheets.Item("Party Roles")
$wsPartyRoles.Range("F11:F12").Select()

# ---------------------------------------------------------------------------
# 6. Re-activate "SAM CPHHolding" as the selected/visible sheet and set the
#    workbook window size/position.
# ---------------------------------------------------------------------------
$wsHolding.Activate()
$excel.ActiveWindow.WindowState = -4143
$excel.Left = 12280
$excel.Top = 8160
$excel.Width = 28110
$excel.Height = 12210
